$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 269.57144
$ws.Range("I9").Value = 188.45454
$ws.Range("K9").Value = 188.45454
$ws.Range("M9").Value = -19.45454000000001

$ws.Range("H58").Value = 281.66666
$ws.Range("I58").Value = 281.66666
$ws.Range("K58").Value = 844.9999799999999
$ws.Range("M58").Value = -694.9999799999999

$ws.Range("H76").Value = 5009.25
$ws.Range("I76").Value = 4052.75
$ws.Range("K76").Value = 4052.75
$ws.Range("M76").Value = -3737.75

$ws.Range("H79").Value = 5009.25
$ws.Range("I79").Value = 4052.75
$ws.Range("K79").Value = 4052.75
$ws.Range("M79").Value = -2960.75

$ws.Range("H92").Value = 1275.7142
$ws.Range("I92").Value = 1312.75
$ws.Range("J92").Value = 1053.5
$ws.Range("K92").Value = 1312.75
$ws.Range("L92").Value = 1053.5
$ws.Range("M92").Value = -64.75
$ws.Range("N92").Value = -3549.5

$ws.Range("H96").Value = 791.8570999999999
$ws.Range("I96").Value = 840.5
$ws.Range("J96").Value = 500
$ws.Range("K96").Value = 2521.5
$ws.Range("L96").Value = 1500
$ws.Range("M96").Value = -1148.5
$ws.Range("N96").Value = -4246

$ws.Range("H107").Value = 167166670
$ws.Range("I107").Value = 333333340
$ws.Range("K107").Value = 333333340
$ws.Range("M107").Value = -333331420

$ws.Range("H132").Value = 2479.157
$ws.Range("I132").Value = 2534.5217
$ws.Range("J132").Value = 1969.8
$ws.Range("K132").Value = 7603.5651
$ws.Range("L132").Value = 5909.4
$ws.Range("M132").Value = -5073.5651
$ws.Range("N132").Value = -10969.4

$ws.Range("H138").Value = 3719.3774
$ws.Range("I138").Value = 3439.2307
$ws.Range("K138").Value = 10317.6921
$ws.Range("M138").Value = -5177.6921

$ws.Range("H141").Value = 2668.55
$ws.Range("I141").Value = 2660.5293
$ws.Range("K141").Value = 7981.5879
$ws.Range("M141").Value = -2801.5879

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8947.037
$ws.Range("I32").Value = 7264.634
$ws.Range("K32").Value = 7264.634
$ws.Range("M32").Value = -6977.634

$ws.Range("H45").Value = 11909973
$ws.Range("J45").Value = 7179.5
$ws.Range("L45").Value = 7179.5
$ws.Range("N45").Value = -7933.5

$ws.Range("H61").Value = 7143.029
$ws.Range("I61").Value = 7695.037
$ws.Range("K61").Value = 7695.037
$ws.Range("M61").Value = -7483.037

$ws.Range("H132").Value = 37082.55
$ws.Range("I132").Value = 8464.429
$ws.Range("J132").Value = 63792.8
$ws.Range("K132").Value = 25393.287
$ws.Range("L132").Value = 191378.4
$ws.Range("M132").Value = -22863.287
$ws.Range("N132").Value = -196438.4

$ws.Range("H136").Value = 7143.029
$ws.Range("I136").Value = 7695.037
$ws.Range("K136").Value = 23085.111
$ws.Range("M136").Value = -20535.111

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4810392.5
$ws.Range("I105").Value = 5684646.5
$ws.Range("K105").Value = 5684646.5
$ws.Range("M105").Value = -5682899.5

$ws.Range("H134").Value = 19934.412
$ws.Range("I134").Value = 20172
$ws.Range("J134").Value = 19498.834
$ws.Range("K134").Value = 60516
$ws.Range("L134").Value = 58496.50199999999
$ws.Range("M134").Value = -57981
$ws.Range("N134").Value = -63566.50199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250

$ws.Range("H58").Value = 7113.5713
$ws.Range("I58").Value = 8232.666999999999
$ws.Range("K58").Value = 8232.666999999999
$ws.Range("M58").Value = -8029.666999999999

$ws.Range("H74").Value = 61570
$ws.Range("J74").Value = 93140
$ws.Range("L74").Value = 93140
$ws.Range("N74").Value = -94888

$ws.Range("H77").Value = 61570
$ws.Range("J77").Value = 93140
$ws.Range("L77").Value = 279420
$ws.Range("N77").Value = -288156

$ws.Range("H105").Value = 605.75
$ws.Range("I105").Value = 524.4545000000001
$ws.Range("K105").Value = 524.4545000000001
$ws.Range("M105").Value = 1222.5455

$ws.Range("H132").Value = 67860.66
$ws.Range("I132").Value = 44096.707
$ws.Range("K132").Value = 132290.121
$ws.Range("M132").Value = -129760.121

$ws.Range("H134").Value = 8799.032999999999
$ws.Range("J134").Value = 18112.715
$ws.Range("L134").Value = 54338.145
$ws.Range("N134").Value = -59408.145

$ws.Range("H136").Value = 7113.5713
$ws.Range("I136").Value = 8232.666999999999
$ws.Range("K136").Value = 24698.001
$ws.Range("M136").Value = -22148.001

$ws.Range("H137").Value = 126349.5
$ws.Range("J137").Value = 126349.5
$ws.Range("L137").Value = 126349.5
$ws.Range("N137").Value = -136549.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 360000100
$ws.Range("I32").Value = 250000130
$ws.Range("J32").Value = 800000000
$ws.Range("K32").Value = 750000390
$ws.Range("L32").Value = 2400000000
$ws.Range("M32").Value = -750000107
$ws.Range("N32").Value = -2400000566

$ws.Range("H46").Value = 114436.57
$ws.Range("I46").Value = 335491.7
$ws.Range("J46").Value = 3909
$ws.Range("K46").Value = 1006475.1
$ws.Range("L46").Value = 11727
$ws.Range("M46").Value = -1006384.1
$ws.Range("N46").Value = -11909

$ws.Range("H87").Value = 16826
$ws.Range("I87").Value = 10478
$ws.Range("K87").Value = 31434
$ws.Range("M87").Value = -30186

$ws.Range("H90").Value = 16826
$ws.Range("I90").Value = 10478
$ws.Range("K90").Value = 94302
$ws.Range("M90").Value = -88062

$ws.Range("H129").Value = 1177822.9
$ws.Range("J129").Value = 1380.4546
$ws.Range("L129").Value = 4141.3638
$ws.Range("N129").Value = -14141.3638

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 8835.41
$ws.Range("I132").Value = 6940.7666
$ws.Range("J132").Value = 15150.889
$ws.Range("K132").Value = 20822.2998
$ws.Range("L132").Value = 45452.667
$ws.Range("M132").Value = -18292.2998
$ws.Range("N132").Value = -50512.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 90658.10000000001
$ws.Range("I22").Value = 222896.5
$ws.Range("K22").Value = 222896.5
$ws.Range("M22").Value = -222601.5

$ws.Range("H27").Value = 90658.10000000001
$ws.Range("I27").Value = 222896.5
$ws.Range("K27").Value = 222896.5
$ws.Range("M27").Value = -222789.5

$ws.Range("H62").Value = 30000
$ws.Range("J62").Value = 30000
$ws.Range("L62").Value = 30000
$ws.Range("N62").Value = -31248

$ws.Range("H65").Value = 30000
$ws.Range("J65").Value = 30000
$ws.Range("L65").Value = 90000
$ws.Range("N65").Value = -96240

$ws.Range("H132").Value = 16121.074
$ws.Range("I132").Value = 18094.783
$ws.Range("J132").Value = 4772.25
$ws.Range("K132").Value = 54284.349
$ws.Range("L132").Value = 14316.75
$ws.Range("M132").Value = -51754.349
$ws.Range("N132").Value = -19376.75

$ws.Range("H136").Value = 91849.61
$ws.Range("I136").Value = 156351
$ws.Range("J136").Value = 7997.8
$ws.Range("K136").Value = 469053
$ws.Range("L136").Value = 23993.4
$ws.Range("M136").Value = -466503
$ws.Range("N136").Value = -29093.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 21000
$ws.Range("J25").Value = 21000
$ws.Range("L25").Value = 21000
$ws.Range("N25").Value = -21586

$ws.Range("H117").Value = 24900
$ws.Range("J117").Value = 24900
$ws.Range("L117").Value = 24900
$ws.Range("N117").Value = -34078

$ws.Range("H132").Value = 14448466
$ws.Range("I132").Value = 15631182
$ws.Range("J132").Value = 1832831.4
$ws.Range("K132").Value = 46893546
$ws.Range("L132").Value = 5498494.199999999
$ws.Range("M132").Value = -46891016
$ws.Range("N132").Value = -5503554.199999999

$ws.Range("H136").Value = 5111.45
$ws.Range("I136").Value = 5075.846
$ws.Range("J136").Value = 6500
$ws.Range("K136").Value = 15227.538
$ws.Range("L136").Value = 19500
$ws.Range("M136").Value = -12677.538
$ws.Range("N136").Value = -24600
